$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.961.31'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.819.43'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.04'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4654'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3663'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07362'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8726'
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.29'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('D12').Value = '1.832.95'
$ws.Range('E12').Value = '  +3.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.397'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07115'
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.512'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.44'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008703'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.67'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').Value = '26.984.72'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('E22').Value = '  -0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.59'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').Value = '2.050.64'
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.02'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.38'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.144'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.248'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.56'
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08901'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7598'
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.165'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.504'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.905'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.089'
$ws.Range('E37').Value = '  -1.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05293'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.975'
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.191'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.342'
$ws.Range('E43').Value = '  -4.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1658'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.454'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4862'
$ws.Range('E46').Value = '  -2.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.42'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.667'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.51'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06297'
$ws.Range('E51').Value = '  +0.03%  '
